$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column B to fit the new longer text
$ws.Columns.Item(2).ColumnWidth = 28.44140625

# Update row 8 "expTitle" value cell to hold the new expected title text,
# re-using the same style as the other expected-title cells (D2/D5)
$ws.Range("B8").Value = "Adactin.com - New User Registration"
$ws.Range("B8").Style = $ws.Range("D2").Style

# Add the new rows for the 4th test case (TC-004 - user registration failure)
$ws.Range("A9").Value = " "

$ws.Range("A10").Value = "TC-004"
$ws.Range("B10").Value = "errorMessage"

$ws.Range("A11").Value = "TC-004"
$ws.Range("B11").Value = "Enter Username"
$ws.Range("B11").Style = $ws.Range("D2").Style

$ws.Range("B11").Select()
